$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -3
    3  = -1
    4  = 1
    5  = 1
    6  = -1
    7  = -1
    10 = 1
    11 = -2
    12 = 5
    13 = 1
    14 = 2
    15 = -1
    16 = 5
    17 = -1
    18 = 3
    19 = 4
    20 = 2
    22 = -1
    23 = -2
    25 = -1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
